$wb = $excel.ActiveWorkbook

# --- Sheet "produtos": append rows 87-89 ---
$produtos = $wb.Worksheets.Item("produtos")

$produtos.Range("A87").Value = 86
$produtos.Range("B87").NumberFormat = "@"
$produtos.Range("B87").Value = "333"
$produtos.Range("D87").Value = "kkk"
$produtos.Range("E87").Value = 33

$produtos.Range("A88").Value = 87
$produtos.Range("B88").Value = "feijão"
$produtos.Range("D88").Value = "KG"
$produtos.Range("E88").Value = 0

$produtos.Range("A89").Value = 88
$produtos.Range("B89").Value = "feijão branco"
$produtos.Range("D89").Value = "KG"
$produtos.Range("E89").Value = 0

# --- Sheet "movimentos": append rows 6-10 ---
$movimentos = $wb.Worksheets.Item("movimentos")

$movimentos.Range("A6").Value = 5
$movimentos.Range("B6").Value = 4
$movimentos.Range("C6").Value = "SAIDA"
$movimentos.Range("D6").Value = 22
$movimentos.Range("E6").Value = "2025-12-04 17:02:35"

$movimentos.Range("A7").Value = 6
$movimentos.Range("B7").Value = 4
$movimentos.Range("C7").Value = "SAIDA"
$movimentos.Range("D7").Value = 2
$movimentos.Range("E7").Value = "2025-12-04 17:02:35"

$movimentos.Range("A8").Value = 7
$movimentos.Range("B8").Value = 9
$movimentos.Range("C8").Value = "SAIDA"
$movimentos.Range("D8").Value = 2
$movimentos.Range("E8").Value = "2025-12-04 17:02:35"

$movimentos.Range("A9").Value = 8
$movimentos.Range("B9").Value = 1
$movimentos.Range("C9").Value = "ENTRADA"
$movimentos.Range("D9").Value = 1
$movimentos.Range("E9").Value = "2025-12-04 17:10:53"

$movimentos.Range("A10").Value = 9
$movimentos.Range("B10").Value = 7
$movimentos.Range("C10").Value = "ENTRADA"
$movimentos.Range("D10").Value = 111
$movimentos.Range("E10").Value = "2025-12-04 17:25:00"
